$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'327.27"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'-1.29%"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'45.15"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'2.34%"
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.Value = "'5.602"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'-2.09%"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'0.08112"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'-2.36%"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'8.718"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'-0.94%"
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.Value = "'4.345"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'-3.43%"
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.Value = "'1.898"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'-3.77%"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'-6.27%"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.Value = "'0.9545"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'3.38%"
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.Value = "'0.1182"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'-5.08%"
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = "'0.1909"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'-1.83%"
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.Value = "'0.1002"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'5.49%"
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "'0.04158"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'5.09%"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'0.1065"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'-0.22%"
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.Value = "'0.001275"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'-2.93%"
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.Value = "'0.005918"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'-1.99%"
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.Value = "'3.594"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'2.34%"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.Value = "'0.3486"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'-0.68%"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = "'8.605"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'-5.67%"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'0.1377"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'0.25%"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'0.88%"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.Value = "'0.04281"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'-3.25%"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.Value = "'0.001242"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.Value = "'0.004563"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'2.71%"
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.Value = "'0.0001236"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'3.53%"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.Value = "'0.0004008"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'0.20%"
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = "'0.02683"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'-4.68%"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'0.05635"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'0.34%"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'25.20%"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.Value = "'0.007713"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'-2.80%"
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = "'0.1399"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'-1.92%"
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.Value = "'0.002064"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'-2.04%"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'0.008695"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'-12.68%"
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.Value = "'0.00007108"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'-3.58%"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.Value = "'0.00000000754"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'0.20%"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'0.003504"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'-2.89%"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'0.002281"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'-0.13%"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'0.00002110"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'0.20%"
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.Value = "'0.0002010"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'0.20%"
$c.Style = "Normal"
